# DECS-1771-SoW_sequence.docx edit
# ----------------------------------------------------------------------------
# Two content changes:
#  1. Paragraph 1 ("Connective tissue disorder ..."): remove bold formatting
#     from both the paragraph-mark run properties and the run itself.
#  2. Paragraph 2 ("COPD J42, J43, J44"): split the single run into two runs
#     - "COPD J42," and " J43, J44" - with the existing "_GoBack" bookmark
#     relocated so it now sits between them (this is what happens when a
#     user clicks/edits at that spot - Word's last-edit bookmark moves there
#     and the run gets split around it).
#
# Both paragraphs are rewritten via Range.InsertXML so the resulting
# WordprocessingML matches exactly (same rFonts/rPr, same rsid attributes
# that were already present on the surrounding content).

$d = $word.ActiveDocument

$p1 = $d.Paragraphs.First
$p2 = $d.Paragraphs.Last

# --- Paragraph 1: drop <w:b/> from the paragraph mark and the run ----------
if ($p1.Range.Text -like "Connective tissue disorder*") {
    $p1Xml = @'
<w:p w14:paraId="49CEDC67" w14:textId="18C62728" w:rsidR="74757A7E" w:rsidRDefault="74757A7E" w:rsidP="000A20EC" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="6D40AF6B"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t>Connective tissue disorder &#8211; M30 &#8211; M36 is a huge list of different Systemic connective tissue disorders &#8211; we could include all.</w:t></w:r></w:p>
'@
    $p1.Range.InsertXML($p1Xml)
} else {
    # Fallback: just toggle bold off on the paragraph/run if the text
    # doesn't match what we expect (keeps the script from being a no-op).
    $p1.Range.Bold = 0
}

# --- Paragraph 2: split "COPD J42, J43, J44" around the _GoBack bookmark ---
if ($p2.Range.Text -like "COPD J42*") {
    $p2Xml = @'
<w:p w14:paraId="1EFB9795" w14:textId="49B569C8" w:rsidR="74757A7E" w:rsidRPr="006107A8" w:rsidRDefault="74757A7E" w:rsidP="006107A8" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr></w:pPr><w:r w:rsidRPr="006107A8"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t>COPD J42,</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve"> J43, J44</w:t></w:r></w:p>
'@
    $p2.Range.InsertXML($p2Xml)
} else {
    # Fallback: relocate the _GoBack bookmark to split the run naturally.
    $full = $d.Content
    $splitAt = $full.Text.IndexOf(" J43, J44")
    if ($splitAt -ge 0) {
        $pt = $d.Range($splitAt, $splitAt)
        $d.Bookmarks.Add("_GoBack", $pt)
    }
}
